# Fix Training Data Issue (#48)
#
# The "Date" column (BF) held strings built as "<game-date>-<season>"
# (e.g. "6-20-2011-12") instead of a real date string. Data was taken
# from 1 day off due to the way NBA stats were shown, so normalize the
# column to ISO "YYYY-MM-DD" text (e.g. "2012-06-20").
#
# NumberFormat is forced to text ("@") right before the write so Excel
# doesn't silently reinterpret the ISO-looking string as a date serial;
# the cell formatting is cleared again immediately after so each cell is
# left exactly as it was found (default/no explicit style), matching the
# rest of the untouched column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-20-2011-12"
$newValue = "2012-06-20"
$dateCol = "BF"

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("$dateCol$r")
    $current = $cell.Value()
    if ($current -eq $oldValue) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    }
}
